$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the stored number-format codes used by style index 2 (date/time -> date only).
#    Style 2 is applied to column B (date) and, in the original file, column H (break time).
$ws.Range('B2:B20').NumberFormat = 'yyyy-mm-dd'

# 2) Re-populate rows 2-20 with the corrected sample (src) data.
#    Column H switches from a date/time serial value to a plain text time-of-day value,
#    matching columns F and G, so its number format is reset to General first.
$ws.Range('H2:H20').NumberFormat = 'General'

$ws.Range('A2').Value = 'Introduction'
$ws.Range('B2').Value = 44562
$ws.Range('C2').Value = 9
$ws.Range('D2').Value = 9
$ws.Range('E2').Value = 1
$ws.Range('F2').Value = '13:00:00'
$ws.Range('G2').Value = '13:45:00'
$ws.Range('H2').Value = '14:00:00'

$ws.Range('A3').Value = 'Kubernetes Overview'
$ws.Range('B3').Value = 44562
$ws.Range('C3').Value = 22
$ws.Range('D3').Value = 31
$ws.Range('E3').Value = 1
$ws.Range('F3').Value = '13:00:00'
$ws.Range('G3').Value = '13:45:00'
$ws.Range('H3').Value = '14:00:00'

$ws.Range('A4').Value = 'Setup Kubernetes'
$ws.Range('B4').Value = 44562
$ws.Range('C4').Value = 14
$ws.Range('D4').Value = 45
$ws.Range('E4').Value = 1
$ws.Range('F4').Value = '13:00:00'
$ws.Range('G4').Value = '13:45:00'
$ws.Range('H4').Value = '14:00:00'

$ws.Range('A5').Value = 'Setup Kubernetes'
$ws.Range('B5').Value = 44562
$ws.Range('C5').Value = 6
$ws.Range('D5').Value = 51
$ws.Range('E5').Value = 2
$ws.Range('F5').Value = '14:00:00'
$ws.Range('G5').Value = '14:45:00'
$ws.Range('H5').Value = '15:00:00'

$ws.Range('A6').Value = 'Kubernetes Concepts'
$ws.Range('B6').Value = 44562
$ws.Range('C6').Value = 13
$ws.Range('D6').Value = 64
$ws.Range('E6').Value = 2
$ws.Range('F6').Value = '14:00:00'
$ws.Range('G6').Value = '14:45:00'
$ws.Range('H6').Value = '15:00:00'

$ws.Range('A7').Value = 'YAML Introduction'
$ws.Range('B7').Value = 44562
$ws.Range('C7').Value = 8
$ws.Range('D7').Value = 72
$ws.Range('E7').Value = 2
$ws.Range('F7').Value = '14:00:00'
$ws.Range('G7').Value = '14:45:00'
$ws.Range('H7').Value = '15:00:00'

$ws.Range('A8').Value = 'Kubernets Concepts - PODs, ReplicaSets, Deployments'
$ws.Range('B8').Value = 44562
$ws.Range('C8').Value = 18
$ws.Range('D8').Value = 90
$ws.Range('E8').Value = 2
$ws.Range('F8').Value = '14:00:00'
$ws.Range('G8').Value = '14:45:00'
$ws.Range('H8').Value = '15:00:00'

$ws.Range('A9').Value = 'Kubernets Concepts - PODs, ReplicaSets, Deployments'
$ws.Range('B9').Value = 44562
$ws.Range('C9').Value = 30
$ws.Range('D9').Value = 120
$ws.Range('E9').Value = 3
$ws.Range('F9').Value = '15:00:00'
$ws.Range('G9').Value = '15:45:00'
$ws.Range('H9').Value = '16:00:00'

$ws.Range('A10').Value = 'Kubernets Concepts - PODs, ReplicaSets, Deployments'
$ws.Range('B10').Value = 44563
$ws.Range('C10').Value = 15
$ws.Range('D10').Value = 135
$ws.Range('E10').Value = 3
$ws.Range('F10').Value = '13:00:00'
$ws.Range('G10').Value = '13:45:00'
$ws.Range('H10').Value = '14:00:00'

$ws.Range('A11').Value = 'Kubernets Concepts - PODs, ReplicaSets, Deployments'
$ws.Range('B11').Value = 44563
$ws.Range('C11').Value = 45
$ws.Range('D11').Value = 180
$ws.Range('E11').Value = 4
$ws.Range('F11').Value = '14:00:00'
$ws.Range('G11').Value = '14:45:00'
$ws.Range('H11').Value = '15:00:00'

$ws.Range('A12').Value = 'Kubernets Concepts - PODs, ReplicaSets, Deployments'
$ws.Range('B12').Value = 44563
$ws.Range('C12').Value = 20
$ws.Range('D12').Value = 200
$ws.Range('E12').Value = 5
$ws.Range('F12').Value = '15:00:00'
$ws.Range('G12').Value = '15:45:00'
$ws.Range('H12').Value = '16:00:00'

$ws.Range('A13').Value = 'Networking in Kubernetes'
$ws.Range('B13').Value = 44563
$ws.Range('C13').Value = 5
$ws.Range('D13').Value = 205
$ws.Range('E13').Value = 5
$ws.Range('F13').Value = '15:00:00'
$ws.Range('G13').Value = '15:45:00'
$ws.Range('H13').Value = '16:00:00'

$ws.Range('A14').Value = 'Services'
$ws.Range('B14').Value = 44563
$ws.Range('C14').Value = 20
$ws.Range('D14').Value = 225
$ws.Range('E14').Value = 5
$ws.Range('F14').Value = '15:00:00'
$ws.Range('G14').Value = '15:45:00'
$ws.Range('H14').Value = '16:00:00'

$ws.Range('A15').Value = 'Services'
$ws.Range('B15').Value = 44563
$ws.Range('C15').Value = 4
$ws.Range('D15').Value = 229
$ws.Range('E15').Value = 6
$ws.Range('F15').Value = '16:00:00'
$ws.Range('G15').Value = '16:45:00'
$ws.Range('H15').Value = '17:00:00'

$ws.Range('A16').Value = 'Microservices Architechture'
$ws.Range('B16').Value = 44563
$ws.Range('C16').Value = 11
$ws.Range('D16').Value = 240
$ws.Range('E16').Value = 6
$ws.Range('F16').Value = '16:00:00'
$ws.Range('G16').Value = '16:45:00'
$ws.Range('H16').Value = '17:00:00'

$ws.Range('A17').Value = 'Microservices Architechture'
$ws.Range('B17').Value = 44564
$ws.Range('C17').Value = 30
$ws.Range('D17').Value = 270
$ws.Range('E17').Value = 6
$ws.Range('F17').Value = '13:00:00'
$ws.Range('G17').Value = '13:45:00'
$ws.Range('H17').Value = '14:00:00'

$ws.Range('A18').Value = 'Microservices Architechture'
$ws.Range('B18').Value = 44564
$ws.Range('C18').Value = 4
$ws.Range('D18').Value = 274
$ws.Range('E18').Value = 7
$ws.Range('F18').Value = '14:00:00'
$ws.Range('G18').Value = '14:45:00'
$ws.Range('H18').Value = '15:00:00'

$ws.Range('A19').Value = 'Kubernetes on the Cloud'
$ws.Range('B19').Value = 44564
$ws.Range('C19').Value = 26
$ws.Range('D19').Value = 300
$ws.Range('E19').Value = 7
$ws.Range('F19').Value = '14:00:00'
$ws.Range('G19').Value = '14:45:00'
$ws.Range('H19').Value = '15:00:00'

$ws.Range('A20').Value = 'Conclusion'
$ws.Range('B20').Value = 44564
$ws.Range('C20').Value = 2
$ws.Range('D20').Value = 302
$ws.Range('E20').Value = 7
$ws.Range('F20').Value = '14:00:00'
$ws.Range('G20').Value = '14:45:00'
$ws.Range('H20').Value = '15:00:00'

